$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1767.6586
$ws.Range("J112").Value = 1774.975
$ws.Range("L112").Value = 5324.924999999999
$ws.Range("N112").Value = -7540.924999999999

$ws.Range("H132").Value = 23353.348
$ws.Range("I132").Value = 1656.4
$ws.Range("J132").Value = 167999.67
$ws.Range("K132").Value = 4969.200000000001
$ws.Range("L132").Value = 503999.01
$ws.Range("M132").Value = -2439.200000000001
$ws.Range("N132").Value = -509059.01

$ws.Range("H137").Value = 2468.65
$ws.Range("I137").Value = 1757.875
$ws.Range("J137").Value = 5311.75
$ws.Range("K137").Value = 5273.625
$ws.Range("L137").Value = 15935.25
$ws.Range("M137").Value = -2723.625
$ws.Range("N137").Value = -21035.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11246.258
$ws.Range("I2").Value = 13339.76
$ws.Range("K2").Value = 13339.76
$ws.Range("M2").Value = -13226.76

$ws.Range("H32").Value = 9833.361000000001
$ws.Range("I32").Value = 7475.271
$ws.Range("K32").Value = 7475.271
$ws.Range("M32").Value = -7188.271

$ws.Range("H45").Value = 4457.75
$ws.Range("I45").Value = 3081
$ws.Range("J45").Value = 6752.3335
$ws.Range("K45").Value = 3081
$ws.Range("L45").Value = 6752.3335
$ws.Range("M45").Value = -2704
$ws.Range("N45").Value = -7506.3335

$ws.Range("H61").Value = 4985.1924
$ws.Range("I61").Value = 3546.4167
$ws.Range("K61").Value = 3546.4167
$ws.Range("M61").Value = -3334.4167

$ws.Range("H74").Value = 2022.5883
$ws.Range("I74").Value = 1798.0333
$ws.Range("K74").Value = 1798.0333
$ws.Range("M74").Value = -924.0333000000001

$ws.Range("H77").Value = 2022.5883
$ws.Range("I77").Value = 1798.0333
$ws.Range("K77").Value = 8990.166499999999
$ws.Range("M77").Value = -4622.166499999999

$ws.Range("H116").Value = 11246.258
$ws.Range("I116").Value = 13339.76
$ws.Range("K116").Value = 13339.76
$ws.Range("M116").Value = -11045.76

$ws.Range("H122").Value = 5819.073
$ws.Range("I122").Value = 5107.8076
$ws.Range("J122").Value = 7051.933
$ws.Range("K122").Value = 15323.4228
$ws.Range("L122").Value = 21155.799
$ws.Range("M122").Value = -12873.4228
$ws.Range("N122").Value = -26055.799

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H136").Value = 4985.1924
$ws.Range("I136").Value = 3546.4167
$ws.Range("K136").Value = 10639.2501
$ws.Range("M136").Value = -8089.250100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11246.258
$ws.Range("I3").Value = 13339.76
$ws.Range("K3").Value = 13339.76
$ws.Range("M3").Value = -13225.76

$ws.Range("H99").Value = 34820.75
$ws.Range("I99").Value = 37622.637
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 37622.637
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = -36124.637
$ws.Range("N99").Value = -6996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

$ws.Range("H99").Value = 18427266
$ws.Range("J99").Value = 20841302
$ws.Range("L99").Value = 20841302
$ws.Range("N99").Value = -20844298

$ws.Range("H126").Value = 18427266
$ws.Range("J126").Value = 20841302
$ws.Range("L126").Value = 62523906
$ws.Range("N126").Value = -62528846

$ws.Range("H132").Value = 4911.7144
$ws.Range("I132").Value = 3376.4
$ws.Range("K132").Value = 10129.2
$ws.Range("M132").Value = -7599.200000000001

$ws.Range("H141").Value = 88878.75999999999
$ws.Range("J141").Value = 91437.47
$ws.Range("L141").Value = 91437.47
$ws.Range("N141").Value = -101797.47

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 123.51613
$ws.Range("J12").Value = 156.83333
$ws.Range("L12").Value = 470.49999
$ws.Range("N12").Value = -816.49999

$ws.Range("H33").Value = 122.2
$ws.Range("I33").Value = 122.2
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 733.2
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -450.2
$ws.Range("N33").ClearContents()

$ws.Range("H68").Value = 23810252
$ws.Range("I68").Value = 922.5
$ws.Range("K68").Value = 2767.5
$ws.Range("M68").Value = -1956.5

$ws.Range("H71").Value = 23810252
$ws.Range("I71").Value = 922.5
$ws.Range("K71").Value = 8302.5
$ws.Range("M71").Value = -4246.5

$ws.Range("H80").Value = 5106.4443
$ws.Range("I80").Value = 4663.3335
$ws.Range("K80").Value = 13990.0005
$ws.Range("M80").Value = -13054.0005

$ws.Range("H83").Value = 5106.4443
$ws.Range("I83").Value = 4663.3335
$ws.Range("K83").Value = 41970.0015
$ws.Range("M83").Value = -37290.0015

$ws.Range("H101").Value = 9500
$ws.Range("J101").Value = 9500
$ws.Range("L101").Value = 28500
$ws.Range("N101").Value = -33368

$ws.Range("H110").Value = 8099.6665
$ws.Range("I110").Value = 8099.6665
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 24298.9995
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -20208.9995
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 285.35715
$ws.Range("I2").Value = 286.17648
$ws.Range("J2").Value = 284.0909
$ws.Range("K2").Value = 286.17648
$ws.Range("L2").Value = 284.0909
$ws.Range("M2").Value = -173.17648
$ws.Range("N2").Value = -510.0909

$ws.Range("H80").Value = 43141.63
$ws.Range("I80").Value = 63218.61
$ws.Range("K80").Value = 63218.61
$ws.Range("M80").Value = -62220.61

$ws.Range("H83").Value = 43141.63
$ws.Range("I83").Value = 63218.61
$ws.Range("K83").Value = 316093.05
$ws.Range("M83").Value = -311101.05

$ws.Range("H97").Value = 324.5
$ws.Range("I97").Value = 324.5
$ws.Range("K97").Value = 324.5
$ws.Range("M97").Value = 171.5

$ws.Range("H132").Value = 4859.5864
$ws.Range("I132").Value = 4093.2942
$ws.Range("J132").Value = 5945.1665
$ws.Range("K132").Value = 12279.8826
$ws.Range("L132").Value = 17835.4995
$ws.Range("M132").Value = -9749.882599999999
$ws.Range("N132").Value = -22895.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4417.9375
$ws.Range("I7").Value = 2463.8333
$ws.Range("J7").Value = 5590.4
$ws.Range("K7").Value = 2463.8333
$ws.Range("L7").Value = 5590.4
$ws.Range("M7").Value = -2351.8333
$ws.Range("N7").Value = -5814.4

$ws.Range("H45").Value = 24666.666
$ws.Range("I45").Value = 18500
$ws.Range("K45").Value = 18500
$ws.Range("M45").Value = -18093

$ws.Range("H46").Value = 4863.275
$ws.Range("I46").Value = 619.5294
$ws.Range("J46").Value = 7999.9565
$ws.Range("K46").Value = 619.5294
$ws.Range("L46").Value = 7999.9565
$ws.Range("M46").Value = -431.5294
$ws.Range("N46").Value = -8375.9565

$ws.Range("H122").Value = 100005080
$ws.Range("I122").Value = 200002220
$ws.Range("J122").Value = 7950
$ws.Range("K122").Value = 600006660
$ws.Range("L122").Value = 23850
$ws.Range("M122").Value = -600004210
$ws.Range("N122").Value = -28750

$ws.Range("H126").Value = 4417.9375
$ws.Range("I126").Value = 2463.8333
$ws.Range("J126").Value = 5590.4
$ws.Range("K126").Value = 7391.499899999999
$ws.Range("L126").Value = 16771.2
$ws.Range("M126").Value = -4921.499899999999
$ws.Range("N126").Value = -21711.2

$ws.Range("H132").Value = 4243.0713
$ws.Range("I132").Value = 3485.3333
$ws.Range("J132").Value = 5607
$ws.Range("K132").Value = 10455.9999
$ws.Range("L132").Value = 16821
$ws.Range("M132").Value = -7925.999899999999
$ws.Range("N132").Value = -21881

$ws.Range("H136").Value = 5349.1875
$ws.Range("I136").Value = 3710.9048
$ws.Range("J136").Value = 8476.817999999999
$ws.Range("K136").Value = 11132.7144
$ws.Range("L136").Value = 25430.454
$ws.Range("M136").Value = -8582.714399999999
$ws.Range("N136").Value = -30530.454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 25000
$ws.Range("J101").Value = 25000
$ws.Range("L101").Value = 25000
$ws.Range("N101").Value = -31490

$ws.Range("H126").Value = 2209.8125
$ws.Range("I126").Value = 2125.8462
$ws.Range("K126").Value = 6377.5386
$ws.Range("M126").Value = -3907.5386

$ws.Range("H136").Value = 4599.794
$ws.Range("I136").Value = 4025.6843
$ws.Range("K136").Value = 12077.0529
$ws.Range("M136").Value = -9527.0529
